$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update case counts (Estados Unidos, Francia, Canada, India, Uzbekistan, Ruanda, Etiopia) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1353534
$ws.Range("C4").Value = 6225
$ws.Range("E4").Value = 1034026

# Row 9: Francia
$ws.Range("D9").Value = 56217
$ws.Range("E9").Value = 94061
$ws.Range("G9").Value = 70
$ws.Range("H9").Value = 26380

# Row 15: Canada
$ws.Range("B15").Value = 68738
$ws.Range("C15").Value = 1036
$ws.Range("D15").Value = 31902
$ws.Range("E15").Value = 31966
$ws.Range("G15").Value = 177
$ws.Range("H15").Value = 4870

# Row 16: India
$ws.Range("B16").Value = 67138
$ws.Range("C16").Value = 4330
$ws.Range("D16").Value = 20969
$ws.Range("E16").Value = 43957
$ws.Range("G16").Value = 111
$ws.Range("H16").Value = 2212

# Row 75: Uzbekistan
$ws.Range("B75").Value = 2418
$ws.Range("C75").Value = 69
$ws.Range("E75").Value = 527

# Row 137: Ruanda
$ws.Range("B137").Value = 284
$ws.Range("C137").Value = 4
$ws.Range("E137").Value = 144

# Row 141: Etiopia
$ws.Range("B141").Value = 241
$ws.Range("C141").Value = 31
$ws.Range("E141").Value = 137

# --- Reorder Belice / Nueva Caledonia (swap rows 192 and 193) ---
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2

# --- Reorder Butan / Islas Virgenes Britanicas (swap rows 212 and 213) ---
$ws.Range("A212").Value = "Islas Virgenes Britanicas"
$ws.Range("D212").Value = 4
$ws.Range("H212").Value = 1

$ws.Range("A213").Value = "Butan"
$ws.Range("D213").Value = 5
$ws.Range("H213").Value = 0
